# fixes for the agent TearDown & MakePayment Test
#
# Applies the set of fixes captured in the commit: a few test-data string
# corrections (customer names, mediation filenames/versions), a removed
# "separator" column on the Mediation sheet, updated counters, and the
# selection / active-sheet bookkeeping that Excel records as a side effect
# of where the author last clicked while making those edits.

$wb = $excel.ActiveWorkbook

# --- AddingAsset: customer name used for the asset test was edited, cursor
#     left on the edited cell (F1). The string itself ("James4950" ->
#     "Asset Customer") is shared, so just re-writing F1 is enough.
$wsAddingAsset = $wb.Sheets.Item("AddingAsset")
$wsAddingAsset.Range("F1").Value = "Asset Customer"
$wsAddingAsset.Range("F1").Select()
$wsAddingAsset.Activate()

# --- Mediation: drop the old "separator count" column E entirely (shifts
#     F:U left into E:T), fix the asterisk file mask & mediation process
#     name, and bump a couple of numeric test parameters.
$wsMediation = $wb.Sheets.Item("Mediation")
$wsMediation.Columns("E:E").Delete()
$wsMediation.Range("E1").Value = "asterisk"
$wsMediation.Range("J1").Value = 2
$wsMediation.Range("M1").Value = 2
$wsMediation.Range("N1").Value = "Web Data"
$wsMediation.Range("Q1").Value = "Test Mediation2.0"
$wsMediation.Range("R1").Value = 12
$wsMediation.Range("A1").Select()
$wsMediation.Activate()

# --- GenerateInvoice: selection only.
$wsGenerateInvoice = $wb.Sheets.Item("GenerateInvoice")
$wsGenerateInvoice.Range("M1").Select()
$wsGenerateInvoice.Activate()

# --- GeneratePayInvoice: selection only.
$wsGeneratePayInvoice = $wb.Sheets.Item("GeneratePayInvoice")
$wsGeneratePayInvoice.Range("T1").Select()
$wsGeneratePayInvoice.Activate()

# --- MakePayment: customer name refreshed to a new unique test value in
#     both places it is used on the row.
$wsMakePayment = $wb.Sheets.Item("MakePayment")
$wsMakePayment.Range("F1").Value = "James7998"
$wsMakePayment.Range("I1").Value = "James7998"
$wsMakePayment.Range("F1").Select()
$wsMakePayment.Activate()

# --- Reports: selection only.
$wsReports = $wb.Sheets.Item("Reports")
$wsReports.Range("R1").Select()
$wsReports.Activate()

# --- Agent: login/customer type corrected from "Web Data Sanity" to
#     "Web Data", row 1 picked up a taller auto-height, and this sheet
#     becomes the workbook's active tab (last Activate() wins).
$wsAgent = $wb.Sheets.Item("Agent")
$wsAgent.Range("H1").Value = "Web Data"
$wsAgent.Rows.Item(1).RowHeight = 15.65
$wsAgent.Range("J1").Select()
$wsAgent.Activate()

# --- AgentPlugin: selection only.
$wsAgentPlugin = $wb.Sheets.Item("AgentPlugin")
$wsAgentPlugin.Range("D1").Select()
$wsAgentPlugin.Activate()

# --- AgentComProcess: selection only.
$wsAgentComProcess = $wb.Sheets.Item("AgentComProcess")
$wsAgentComProcess.Range("E1").Select()
$wsAgentComProcess.Activate()

# --- TearDown: selection only.
$wsTearDown = $wb.Sheets.Item("TearDown")
$wsTearDown.Range("D1").Select()
$wsTearDown.Activate()

# Final active tab is Agent (workbook activeTab 24 -> 34), matching the
# OrderHierarcy -> Agent tab-selection swap in the diff.
$wsAgent.Activate()
